$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.133.06"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "3.831.73"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'704.11"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").Value = "'171.89"
$ws.Range("E6").Value = "  -1.29%  "
$ws.Range("D7").Value = "3.830.93"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("D11").Value = "'7.39"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("D14").Value = "'36.65"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").Value = "4.478.98"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "3.790.81"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("D17").Value = "71.118.68"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").Value = "'7.23"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "'17.39"
$ws.Range("E20").Value = "  -2.19%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'495.24"
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'10.70"
$ws.Range("E22").Value = "  -3.99%  "
$ws.Range("D23").Value = "'0.736"
$ws.Range("E23").Value = "  +2.66%  "
$ws.Range("D24").Value = "'85.28"
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("E25").Value = "  +0.85%  "
$ws.Range("D26").Value = "'10.63"
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("D27").Value = "'12.11"
$ws.Range("E27").Value = "  -2.05%  "
$ws.Range("E28").Value = "  -2.99%  "
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("D31").Value = "'7.44"
$ws.Range("E31").Value = "  -2.09%  "
$ws.Range("D32").Value = "'2.23"
$ws.Range("E32").Value = "  -3.33%  "
$ws.Range("D33").Value = "'29.40"
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("D34").Value = "'0.176"
$ws.Range("E34").Value = "  -3.56%  "
$ws.Range("D35").Value = "'9.20"
$ws.Range("E35").Value = "  -1.23%  "
$ws.Range("E36").Value = "  +0.71%  "
$ws.Range("D37").Value = "3.795.08"
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("D39").Value = "'2.34"
$ws.Range("E39").Value = "  -2.75%  "
$ws.Range("E40").Value = "  +4.50%  "
$ws.Range("D41").Value = "'5.99"
$ws.Range("E41").Value = "  -0.95%  "
$ws.Range("E42").Value = "  -2.92%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("D46").Value = "'163.86"
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D47").Value = "'428.60"
$ws.Range("E47").Value = "  +4.20%  "
$ws.Range("D48").Value = "'48.89"
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("E51").Value = "  -1.83%  "
